# Auto-generated script applying odds updates to Jogos_do_Dia_Betfair_Back_Lay_2026-01-19 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.51
$ws.Range("G2").Value = 1.57
$ws.Range("H2").Value = 6.6
$ws.Range("I2").Value = 7.8
$ws.Range("P2").Value = 2.3
$ws.Range("R2").Value = 1.52
$ws.Range("T2").Value = 1.8
$ws.Range("U2").Value = 2.06
$ws.Range("V2").Value = 1.15
$ws.Range("W2").Value = 2.74
$ws.Range("Y2").Value = 29
$ws.Range("Z2").Value = 65
$ws.Range("AA2").Value = 230
$ws.Range("AC2").Value = 11.5
$ws.Range("AD2").Value = 27
$ws.Range("AF2").Value = 11
$ws.Range("AI2").Value = 85
$ws.Range("AJ2").Value = 14.5
$ws.Range("AL2").Value = 32
$ws.Range("AN2").Value = 6.8
$ws.Range("AO2").Value = 130
# Row 3
$ws.Range("F3").Value = 5.1
$ws.Range("G3").Value = 6.4
$ws.Range("I3").Value = 2.16
$ws.Range("J3").Value = 2.8
$ws.Range("L3").Value = 1.6
$ws.Range("N3").Value = 2.36
$ws.Range("O3").Value = 1.58
$ws.Range("Q3").Value = 2.74
$ws.Range("S3").Value = 5.4
$ws.Range("T3").Value = 2.24
$ws.Range("U3").Value = 1.65
$ws.Range("V3").Value = 1.87
$ws.Range("AA3").Value = 65
$ws.Range("AB3").Value = 14.5
$ws.Range("AE3").Value = 65
$ws.Range("AG3").Value = 27
$ws.Range("AH3").Value = 65
$ws.Range("AI3").Value = 250
# Row 4
$ws.Range("F4").Value = 2.78
$ws.Range("H4").Value = 2.46
$ws.Range("L4").Value = 1.33
$ws.Range("P4").Value = 2.38
$ws.Range("Q4").Value = 1.67
$ws.Range("R4").Value = 1.52
$ws.Range("U4").Value = 2.46
$ws.Range("W4").Value = 1.5
$ws.Range("X4").Value = 90
$ws.Range("AE4").Value = 980
$ws.Range("AN4").Value = 600
# Row 5
$ws.Range("H5").Value = 7
$ws.Range("N5").Value = 4.2
$ws.Range("P5").Value = 2.1
$ws.Range("T5").Value = 1.95
$ws.Range("U5").Value = 1.98
$ws.Range("X5").Value = 17.5
$ws.Range("Y5").Value = 24
$ws.Range("AE5").Value = 240
$ws.Range("AG5").Value = 9.6
$ws.Range("AH5").Value = 23
$ws.Range("AJ5").Value = 14.5
$ws.Range("AM5").Value = 580
$ws.Range("AO5").Value = 140
# Row 6
$ws.Range("G6").Value = 5.4
$ws.Range("H6").Value = 1.8
$ws.Range("I6").Value = 1.82
$ws.Range("V6").Value = 2.22
$ws.Range("Z6").Value = 10.5
$ws.Range("AO6").Value = 55
# Row 7
$ws.Range("I7").Value = 11.5
$ws.Range("L7").Value = 1.43
$ws.Range("S7").Value = 3.95
$ws.Range("T7").Value = 2.32
$ws.Range("U7").Value = 1.62
# Row 8
$ws.Range("F8").Value = 1.68
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 3.25
$ws.Range("O8").Value = 1.37
$ws.Range("P8").Value = 1.77
$ws.Range("Q8").Value = 2.04
$ws.Range("R8").Value = 1.29
$ws.Range("S8").Value = 3.75
$ws.Range("T8").Value = 1.98
$ws.Range("U8").Value = 1.84
$ws.Range("W8").Value = 2.22
$ws.Range("AH8").Value = 48
# Row 9
$ws.Range("N9").Value = 3.35
# Row 10
$ws.Range("F10").Value = 1.4
$ws.Range("H10").Value = 8.199999999999999
$ws.Range("I10").Value = 8.800000000000001
$ws.Range("N10").Value = 6.4
$ws.Range("P10").Value = 2.78
$ws.Range("Q10").Value = 1.51
$ws.Range("R10").Value = 1.75
$ws.Range("S10").Value = 2.24
$ws.Range("T10").Value = 1.76
$ws.Range("U10").Value = 2.2
$ws.Range("X10").Value = 32
$ws.Range("AD10").Value = 80
$ws.Range("AE10").Value = 510
$ws.Range("AF10").Value = 10.5
$ws.Range("AI10").Value = 85
$ws.Range("AM10").Value = 95
# Row 11
$ws.Range("F11").Value = 2.26
$ws.Range("G11").Value = 2.44
$ws.Range("I11").Value = 3.45
$ws.Range("J11").Value = 3.55
$ws.Range("K11").Value = 3.85
$ws.Range("L11").Value = 1.36
$ws.Range("P11").Value = 1.9
$ws.Range("V11").Value = 1.41
$ws.Range("AM11").Value = 330
# Row 12
$ws.Range("F12").Value = 2.68
$ws.Range("G12").Value = 2.7
$ws.Range("H12").Value = 3.3
$ws.Range("I12").Value = 3.35
$ws.Range("P12").Value = 1.63
$ws.Range("S12").Value = 5
$ws.Range("V12").Value = 1.42
$ws.Range("AE12").Value = 44
$ws.Range("AJ12").Value = 38
$ws.Range("AK12").Value = 34
# Row 13
$ws.Range("F13").Value = 4.8
$ws.Range("G13").Value = 6.2
$ws.Range("H13").Value = 1.77
$ws.Range("I13").Value = 1.88
$ws.Range("J13").Value = 3.4
$ws.Range("K13").Value = 3.95
$ws.Range("O13").Value = 1.38
$ws.Range("P13").Value = 1.59
$ws.Range("Q13").Value = 2.18
$ws.Range("S13").Value = 4.2
$ws.Range("V13").Value = 2.12
$ws.Range("W13").Value = 1.2
$ws.Range("X13").Value = 12
$ws.Range("AB13").Value = 970
$ws.Range("AC13").Value = 14
$ws.Range("AH13").Value = 60
# Row 14
$ws.Range("F14").Value = 2.1
$ws.Range("J14").Value = 3.3
$ws.Range("K14").Value = 3.55
$ws.Range("R14").Value = 1.26
$ws.Range("W14").Value = 1.83
# Row 15
$ws.Range("I15").Value = 2
$ws.Range("K15").Value = 4.1
$ws.Range("Q15").Value = 1.74
$ws.Range("S15").Value = 2.96
$ws.Range("U15").Value = 2.16
$ws.Range("V15").Value = 1.96
$ws.Range("W15").Value = 1.27
# Row 16
$ws.Range("F16").Value = 2.22
$ws.Range("G16").Value = 2.32
$ws.Range("K16").Value = 3.25
$ws.Range("T16").Value = 1.98
$ws.Range("AB16").Value = 14
# Row 17
$ws.Range("H17").Value = 3.25
# Row 18
$ws.Range("F18").Value = 1.93
$ws.Range("H18").Value = 3.45
$ws.Range("J18").Value = 4.4
$ws.Range("P18").Value = 3.4
$ws.Range("R18").Value = 1.92
$ws.Range("S18").Value = 1.76
$ws.Range("U18").Value = 3.05
# Row 19
$ws.Range("F19").Value = 3.7
$ws.Range("P19").Value = 2.76
$ws.Range("S19").Value = 2.18
$ws.Range("V19").Value = 1.99
# Row 20
$ws.Range("L20").Value = 1.45
$ws.Range("AG20").Value = 11
# Row 21
$ws.Range("F21").Value = 2.4
$ws.Range("G21").Value = 2.48
$ws.Range("H21").Value = 3.55
$ws.Range("K21").Value = 3.2
$ws.Range("N21").Value = 2.7
$ws.Range("W21").Value = 1.67
$ws.Range("AA21").Value = 80
$ws.Range("AG21").Value = 12.5
# Row 22
$ws.Range("N22").Value = 3.3
$ws.Range("P22").Value = 1.84
# Row 23
$ws.Range("F23").Value = 2.98
$ws.Range("P23").Value = 1.65
$ws.Range("X23").Value = 9.199999999999999
$ws.Range("Y23").Value = 9
$ws.Range("Z23").Value = 18
$ws.Range("AD23").Value = 13
$ws.Range("AE23").Value = 38
$ws.Range("AF23").Value = 18.5
# Row 24
$ws.Range("L24").Value = 1.28
$ws.Range("P24").Value = 2.8
$ws.Range("Q24").Value = 1.54
$ws.Range("T24").Value = 1.55
$ws.Range("AB24").Value = 14.5
$ws.Range("AL24").Value = 24
# Row 25
$ws.Range("Q25").Value = 2.22
$ws.Range("AA25").Value = 65
$ws.Range("AE25").Value = 44
$ws.Range("AK25").Value = 26
$ws.Range("AL25").Value = 44
# Row 26
$ws.Range("F26").Value = 2.94
# Row 27
$ws.Range("G27").Value = 2.46
$ws.Range("I27").Value = 3.85
$ws.Range("O27").Value = 1.44
$ws.Range("Q27").Value = 2.3
$ws.Range("W27").Value = 1.68
